$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 192, shifting existing rows 192..287 down to 193..288.
$ws.Rows.Item(192).Insert()

# Populate the newly inserted row 192 with the new record's data.
$ws.Cells.Item(192, 1).Value2 = 7
$ws.Cells.Item(192, 2).Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(192, 3).Value2 = 'Ñuble'
$ws.Cells.Item(192, 4).Value2 = 44609
$ws.Cells.Item(192, 5).Value2 = 16
$ws.Cells.Item(192, 6).Value2 = 100114001
$ws.Cells.Item(192, 7).Value2 = 'Papa'
$ws.Cells.Item(192, 8).Value2 = 'Patagonia'
$ws.Cells.Item(192, 9).Value2 = '1a nueva(o)'
$ws.Cells.Item(192, 10).Value2 = 200
$ws.Cells.Item(192, 11).Value2 = 6500
$ws.Cells.Item(192, 12).Value2 = 7000
$ws.Cells.Item(192, 13).Value2 = 6750
$ws.Cells.Item(192, 14).Value2 = '$/saco 25 kilos'
$ws.Cells.Item(192, 15).Value2 = 'Provincia de Diguillín'
$ws.Cells.Item(192, 16).Value2 = 270
$ws.Cells.Item(192, 17).Value2 = 25
$ws.Cells.Item(192, 18).Value2 = 'Hortaliza'
